$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''55.116.38'
$ws.Range("E2").Value = '  +4.24%  '

$ws.Range("D3").Value = '''2.445.63'
$ws.Range("E3").Value = '  +3.87%  '

$ws.Range("E4").Value = '  +0.45%  '

$ws.Range("D5").Value = '''477.75'
$ws.Range("E5").Value = '  +6.92%  '

$ws.Range("D6").Value = '''139.60'
$ws.Range("E6").Value = '  +11.72%  '

$ws.Range("D7").Value = '''0.997'
$ws.Range("E7").Value = '  +0.25%  '

$ws.Range("D8").Value = '''0.503'
$ws.Range("E8").Value = '  +6.53%  '

$ws.Range("D9").Value = '''2.455.48'
$ws.Range("E9").Value = '  +4.94%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '''0.0958'
$ws.Range("E10").Value = '  +5.35%  '

$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").Value = '''5.52'
$ws.Range("E11").Value = '  +6.37%  '

$ws.Range("D12").Value = '''0.324'
$ws.Range("E12").Value = '  +4.91%  '

$ws.Range("D13").Value = '''0.123'
$ws.Range("E13").Value = '  +2.08%  '

$ws.Range("D14").Value = '''2.859.89'
$ws.Range("E14").Value = '  +4.79%  '

$ws.Range("D15").Value = '''55.167.81'
$ws.Range("E15").Value = '  +4.80%  '

$ws.Range("D16").Value = '''20.41'
$ws.Range("E16").Value = '  +6.46%  '

$ws.Range("D17").Value = '''0.0000133'
$ws.Range("E17").Value = '  +9.21%  '

$ws.Range("D18").Value = '''2.448.67'
$ws.Range("E18").Value = '  +4.72%  '

$ws.Range("D19").Value = '''4.33'
$ws.Range("E19").Value = '  +5.61%  '

$ws.Range("D20").Value = '''9.90'
$ws.Range("E20").Value = '  +9.74%  '

$ws.Range("D21").Value = '''312.73'
$ws.Range("E21").Value = '  +3.76%  '

$ws.Range("D22").Value = '''0.994'
$ws.Range("E22").Value = '  -1.65%  '

$ws.Range("D23").Value = '''5.68'
$ws.Range("E23").Value = '  +7.91%  '

$ws.Range("D24").Value = '''57.07'
$ws.Range("E24").Value = '  +5.00%  '

$ws.Range("D25").Value = '''0.997'
$ws.Range("E25").Value = '  +0.19%  '

$ws.Range("D26").Value = '''0.400'
$ws.Range("E26").Value = '  +6.37%  '

$ws.Range("D27").Value = '''0.161'
$ws.Range("E27").Value = '  +16.75%  '

$ws.Range("D28").Value = '''2.532.40'
$ws.Range("E28").Value = '  +6.75%  '

$ws.Range("D29").Value = '''7.30'
$ws.Range("E29").Value = '  +4.97%  '

$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '''0.0' + [char]0x2083 + '0765'
$ws.Range("E30").Value = '  +12.53%  '

$ws.Range("B31").Value = 'USDe'
$ws.Range("C31").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D31").Value = '''0.998'
$ws.Range("E31").Value = '  +0.14%  '

$ws.Range("D32").Value = '''148.00'
$ws.Range("E32").Value = '  +0.96%  '

$ws.Range("D33").Value = '''18.12'
$ws.Range("E33").Value = '  +4.93%  '

$ws.Range("D34").Value = '''1.46'
$ws.Range("E34").Value = '  +8.81%  '

$ws.Range("D35").Value = '''5.10'
$ws.Range("E35").Value = '  +5.45%  '

$ws.Range("D36").Value = '''1.12'
$ws.Range("E36").Value = '  +10.64%  '

$ws.Range("D37").Value = '''3.58'
$ws.Range("E37").Value = '  +5.52%  '

$ws.Range("D38").Value = '''0.837'
$ws.Range("E38").Value = '  +3.92%  '

$ws.Range("D39").Value = '''33.63'
$ws.Range("E39").Value = '  +4.51%  '

$ws.Range("D40").Value = '''0.993'
$ws.Range("E40").Value = '  +0.69%  '

$ws.Range("D41").Value = '''3.42'
$ws.Range("E41").Value = '  +6.79%  '

$ws.Range("D42").Value = '''0.599'
$ws.Range("E42").Value = '  +5.78%  '

$ws.Range("D43").Value = '''0.0545'
$ws.Range("E43").Value = '  +6.50%  '

$ws.Range("D44").Value = '''1.28'
$ws.Range("E44").Value = '  +8.01%  '

$ws.Range("D45").Value = '''10.15'
$ws.Range("E45").Value = '  +0.45%  '

$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").Value = '''4.66'
$ws.Range("E46").Value = '  +11.55%  '

$ws.Range("B47").Value = 'Bittensor'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D47").Value = '''255.11'
$ws.Range("E47").Value = '  +26.70%  '

$ws.Range("D48").Value = '''0.0891'
$ws.Range("E48").Value = '  +7.56%  '

$ws.Range("D49").Value = '''0.0221'
$ws.Range("E49").Value = '  +5.29%  '

$ws.Range("D50").Value = '''1.896.03'
$ws.Range("E50").Value = '  -0.96%  '

$ws.Range("D51").Value = '''16.97'
$ws.Range("E51").Value = '  +5.34%  '

